# Apply the commit "from v0.2 to v1.0.1":
#  1. Update the Version cell (D2) from "0.1" to "1.0.1"
#  2. Swap the "step 2" contents between TC2 (row 20) and TC3 (row 28):
#     - Row 20 (TC2, step 2) gets the "realizar a liquidação" content
#       that used to be TC3's step 2.
#     - Row 28 (TC3, step 2) gets the "atribuir/desatribuir" content
#       that used to be TC2's step 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump version string
$ws.Range("D2").Value = "1.0.1"

# 2. Swap the step-2 descriptions / expected results between TC2 and TC3
$oldB20 = $ws.Range("B20").Value2
$oldD20 = $ws.Range("D20").Value2
$oldB28 = $ws.Range("B28").Value2
$oldD28 = $ws.Range("D28").Value2

$ws.Range("B20").Value = $oldB28
$ws.Range("D20").Value = $oldD28
$ws.Range("B28").Value = $oldB20
$ws.Range("D28").Value = $oldD20
